# edit.ps1 - apply "Update countries & provincias Spain" changes
# Source workbook: paises.xlsx (sheet "Pais")
#
# 1) Bump the "last updated" timestamp banner in A1.
# 2) A block of countries got re-sorted/re-inserted into the shared-string
#    table, which (together with refreshed case counts) changes the
#    country name and/or the 7 numeric columns (B:H) on a set of existing
#    rows. We reproduce that final state by writing the new value into
#    every cell that actually differs from the original workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 29 de Marzo de 2020 a las 14:50"

# Row 19: Canada
$ws.Cells.Item(19, 5).Value = 5086
$ws.Cells.Item(19, 7).Value = 1
$ws.Cells.Item(19, 8).Value = 61

# Row 20: Noruega
$ws.Cells.Item(20, 2).Value = 4232
$ws.Cells.Item(20, 3).Value = 217
$ws.Cells.Item(20, 5).Value = 4200

# Row 25: Chequia
$ws.Cells.Item(25, 2).Value = 2697
$ws.Cells.Item(25, 3).Value = 66
$ws.Cells.Item(25, 5).Value = 2673

# Row 32: Polonia
$ws.Cells.Item(32, 1).Value = "Polonia"
$ws.Cells.Item(32, 2).Value = 1771
$ws.Cells.Item(32, 3).Value = 133
$ws.Cells.Item(32, 4).Value = 7
$ws.Cells.Item(32, 5).Value = 1744
$ws.Cells.Item(32, 6).Value = 3
$ws.Cells.Item(32, 7).Value = 2
$ws.Cells.Item(32, 8).Value = 20

# Row 33: Rumania
$ws.Cells.Item(33, 1).Value = "Rumania"
$ws.Cells.Item(33, 2).Value = 1760
$ws.Cells.Item(33, 3).Value = 308
$ws.Cells.Item(33, 4).Value = 169
$ws.Cells.Item(33, 5).Value = 1551
$ws.Cells.Item(33, 6).Value = 34
$ws.Cells.Item(33, 7).Value = 3
$ws.Cells.Item(33, 8).Value = 40

# Row 39: Arabia Saudita
$ws.Cells.Item(39, 1).Value = "Arabia Saudita"
$ws.Cells.Item(39, 2).Value = 1299
$ws.Cells.Item(39, 3).Value = 96
$ws.Cells.Item(39, 4).Value = 66
$ws.Cells.Item(39, 5).Value = 1225
$ws.Cells.Item(39, 6).Value = 6
$ws.Cells.Item(39, 7).Value = 4
$ws.Cells.Item(39, 8).Value = 8

# Row 40: Indonesia
$ws.Cells.Item(40, 1).Value = "Indonesia"
$ws.Cells.Item(40, 2).Value = 1285
$ws.Cells.Item(40, 3).Value = 130
$ws.Cells.Item(40, 4).Value = 64
$ws.Cells.Item(40, 5).Value = 1107
$ws.Cells.Item(40, 6).Value = 0
$ws.Cells.Item(40, 7).Value = 12
$ws.Cells.Item(40, 8).Value = 114

# Row 41: Finlandia
$ws.Cells.Item(41, 1).Value = "Finlandia"
$ws.Cells.Item(41, 2).Value = 1221
$ws.Cells.Item(41, 3).Value = 54
$ws.Cells.Item(41, 4).Value = 10
$ws.Cells.Item(41, 5).Value = 1200
$ws.Cells.Item(41, 6).Value = 32
$ws.Cells.Item(41, 7).Value = 2
$ws.Cells.Item(41, 8).Value = 11

# Row 48: Singapur
$ws.Cells.Item(48, 2).Value = 844
$ws.Cells.Item(48, 3).Value = 42
$ws.Cells.Item(48, 4).Value = 212
$ws.Cells.Item(48, 5).Value = 629

# Row 66: Marruecos
$ws.Cells.Item(66, 1).Value = "Marruecos"
$ws.Cells.Item(66, 2).Value = 450
$ws.Cells.Item(66, 3).Value = 48
$ws.Cells.Item(66, 4).Value = 13
$ws.Cells.Item(66, 5).Value = 411
$ws.Cells.Item(66, 6).Value = 1
$ws.Cells.Item(66, 7).Value = 1
$ws.Cells.Item(66, 8).Value = 26

# Row 67: Libano
$ws.Cells.Item(67, 1).Value = "Libano"
$ws.Cells.Item(67, 2).Value = 438
$ws.Cells.Item(67, 3).Value = 26
$ws.Cells.Item(67, 4).Value = 30
$ws.Cells.Item(67, 5).Value = 398
$ws.Cells.Item(67, 6).Value = 4
$ws.Cells.Item(67, 7).Value = 2
$ws.Cells.Item(67, 8).Value = 10

# Row 68: Lituania
$ws.Cells.Item(68, 1).Value = "Lituania"
$ws.Cells.Item(68, 3).Value = 43
$ws.Cells.Item(68, 4).Value = 1
$ws.Cells.Item(68, 5).Value = 429
$ws.Cells.Item(68, 6).Value = 2
$ws.Cells.Item(68, 7).Value = 0
$ws.Cells.Item(68, 8).Value = 7

# Row 136: Polinesia Francesa
$ws.Cells.Item(136, 1).Value = "Polinesia Francesa"

# Row 137: Uganda
$ws.Cells.Item(137, 1).Value = "Uganda"

# Row 144: Congo
$ws.Cells.Item(144, 1).Value = "Congo"
$ws.Cells.Item(144, 3).Value = 15
$ws.Cells.Item(144, 4).Value = 0
$ws.Cells.Item(144, 5).Value = 19

# Row 145: Etiopia
$ws.Cells.Item(145, 1).Value = "Etiopia"
$ws.Cells.Item(145, 2).Value = 19
$ws.Cells.Item(145, 3).Value = 3
$ws.Cells.Item(145, 4).Value = 1
$ws.Cells.Item(145, 5).Value = 18
$ws.Cells.Item(145, 8).Value = 0

# Row 146: Mali
$ws.Cells.Item(146, 1).Value = "Mali"
$ws.Cells.Item(146, 3).Value = 0

# Row 147: Niger
$ws.Cells.Item(147, 1).Value = "Niger"
$ws.Cells.Item(147, 2).Value = 18
$ws.Cells.Item(147, 3).Value = 8
$ws.Cells.Item(147, 8).Value = 1

# Row 148: Islas Virgenes de los Estados Unidos
$ws.Cells.Item(148, 1).Value = "Islas Virgenes de los Estados Unidos"
$ws.Cells.Item(148, 3).Value = 0
$ws.Cells.Item(148, 4).Value = 0
$ws.Cells.Item(148, 5).Value = 17

# Row 149: Maldivas
$ws.Cells.Item(149, 1).Value = "Maldivas"
$ws.Cells.Item(149, 2).Value = 17
$ws.Cells.Item(149, 3).Value = 1
$ws.Cells.Item(149, 4).Value = 11
$ws.Cells.Item(149, 5).Value = 6

# Row 150: Guinea
$ws.Cells.Item(150, 1).Value = "Guinea"
$ws.Cells.Item(150, 2).Value = 16
$ws.Cells.Item(150, 3).Value = 8
$ws.Cells.Item(150, 5).Value = 16

# Row 151: Nueva Caledonia
$ws.Cells.Item(151, 1).Value = "Nueva Caledonia"
$ws.Cells.Item(151, 3).Value = 0
$ws.Cells.Item(151, 4).Value = 0
$ws.Cells.Item(151, 5).Value = 15

# Row 152: Haiti
$ws.Cells.Item(152, 1).Value = "Haiti"
$ws.Cells.Item(152, 2).Value = 15
$ws.Cells.Item(152, 3).Value = 7
$ws.Cells.Item(152, 4).Value = 1

# Row 153: Republica de Yibuti
$ws.Cells.Item(153, 1).Value = "Republica de Yibuti"
$ws.Cells.Item(153, 4).Value = 0
$ws.Cells.Item(153, 5).Value = 14

# Row 154: Tanzania
$ws.Cells.Item(154, 1).Value = "Tanzania"
$ws.Cells.Item(154, 2).Value = 14
$ws.Cells.Item(154, 4).Value = 1
$ws.Cells.Item(154, 5).Value = 13

# Row 156: Mongolia
$ws.Cells.Item(156, 1).Value = "Mongolia"
$ws.Cells.Item(156, 2).Value = 12
$ws.Cells.Item(156, 5).Value = 12

# Row 157: Dominica
$ws.Cells.Item(157, 1).Value = "Dominica"

# Row 158: San Martin (Parte Francesa)
$ws.Cells.Item(158, 1).Value = "San Martin (Parte Francesa)"
$ws.Cells.Item(158, 3).Value = 0
$ws.Cells.Item(158, 4).Value = 0
$ws.Cells.Item(158, 5).Value = 11

# Row 159: Namibia
$ws.Cells.Item(159, 1).Value = "Namibia"
$ws.Cells.Item(159, 2).Value = 11
$ws.Cells.Item(159, 3).Value = 3
$ws.Cells.Item(159, 4).Value = 2

# Row 160: Bahamas
$ws.Cells.Item(160, 1).Value = "Bahamas"
$ws.Cells.Item(160, 4).Value = 1
$ws.Cells.Item(160, 5).Value = 9

# Row 161: Groenlandia
$ws.Cells.Item(161, 1).Value = "Groenlandia"
$ws.Cells.Item(161, 2).Value = 10
$ws.Cells.Item(161, 4).Value = 2
$ws.Cells.Item(161, 5).Value = 8

# Row 163: Suazilandia
$ws.Cells.Item(163, 1).Value = "Suazilandia"
$ws.Cells.Item(163, 2).Value = 9
$ws.Cells.Item(163, 5).Value = 9

# Row 164: Laos
$ws.Cells.Item(164, 1).Value = "Laos"

# Row 165: Seychelles
$ws.Cells.Item(165, 1).Value = "Seychelles"

# Row 166: Birmania
$ws.Cells.Item(166, 1).Value = "Birmania"

# Row 167: Surinam
$ws.Cells.Item(167, 1).Value = "Surinam"

# Row 168: Mozambique
$ws.Cells.Item(168, 1).Value = "Mozambique"
$ws.Cells.Item(168, 5).Value = 8
$ws.Cells.Item(168, 8).Value = 0

# Row 170: Guyana
$ws.Cells.Item(170, 1).Value = "Guyana"
$ws.Cells.Item(170, 4).Value = 0
$ws.Cells.Item(170, 5).Value = 7

# Row 171: Curazao
$ws.Cells.Item(171, 1).Value = "Curazao"
$ws.Cells.Item(171, 2).Value = 8
$ws.Cells.Item(171, 4).Value = 2
$ws.Cells.Item(171, 5).Value = 5
$ws.Cells.Item(171, 8).Value = 1

# Row 172: Antigua y Barbuda
$ws.Cells.Item(172, 1).Value = "Antigua y Barbuda"
$ws.Cells.Item(172, 5).Value = 7
$ws.Cells.Item(172, 8).Value = 0

# Row 173: Gabon
$ws.Cells.Item(173, 1).Value = "Gabon"

# Row 174: Zimbabue
$ws.Cells.Item(174, 1).Value = "Zimbabue"
$ws.Cells.Item(174, 2).Value = 7
$ws.Cells.Item(174, 8).Value = 1

# Row 175: Benin
$ws.Cells.Item(175, 1).Value = "Benin"

# Row 176: Santa Sede
$ws.Cells.Item(176, 1).Value = "Santa Sede"

# Row 177: Eritrea
$ws.Cells.Item(177, 1).Value = "Eritrea"
$ws.Cells.Item(177, 5).Value = 6
$ws.Cells.Item(177, 8).Value = 0

# Row 178: Cabo Verde
$ws.Cells.Item(178, 1).Value = "Cabo Verde"
$ws.Cells.Item(178, 2).Value = 6
$ws.Cells.Item(178, 8).Value = 1

# Row 180: Montserrat
$ws.Cells.Item(180, 1).Value = "Montserrat"

# Row 182: Fiyi
$ws.Cells.Item(182, 1).Value = "Fiyi"

# Row 183: Siria
$ws.Cells.Item(183, 1).Value = "Siria"

# Row 184: San Bartolome
$ws.Cells.Item(184, 1).Value = "San Bartolome"
$ws.Cells.Item(184, 5).Value = 5
$ws.Cells.Item(184, 8).Value = 0

# Row 185: Sudan
$ws.Cells.Item(185, 1).Value = "Sudan"
$ws.Cells.Item(185, 4).Value = 0
$ws.Cells.Item(185, 8).Value = 1

# Row 186: Nepal
$ws.Cells.Item(186, 1).Value = "Nepal"
$ws.Cells.Item(186, 2).Value = 5
$ws.Cells.Item(186, 4).Value = 1

# Row 189: Santa Lucia
$ws.Cells.Item(189, 1).Value = "Santa Lucia"
$ws.Cells.Item(189, 3).Value = 1
$ws.Cells.Item(189, 4).Value = 1
$ws.Cells.Item(189, 8).Value = 0

# Row 190: Nicaragua
$ws.Cells.Item(190, 1).Value = "Nicaragua"
$ws.Cells.Item(190, 3).Value = 0
$ws.Cells.Item(190, 4).Value = 0
$ws.Cells.Item(190, 8).Value = 1

# Row 191: Somalia
$ws.Cells.Item(191, 1).Value = "Somalia"

# Row 196: Republica del Chad
$ws.Cells.Item(196, 1).Value = "Republica del Chad"
